$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Prepare new header cells (F1:O1) with the same formatting as existing headers ---
$ws.Range("B1").Copy()
$ws.Range("C1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Prepare new species rows (A6:A24) with the same formatting as existing species cells ---
$ws.Range("A5").Copy()
$ws.Range("A6:A24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Set the header row values (A1:O1) ---
$ws.Range("A1").Value = "Species"
$ws.Range("B1").Value = "Mating System"
$ws.Range("C1").Value = "Social Dominance Hierarchy"
$ws.Range("D1").Value = "Territoriality (males)"
$ws.Range("E1").Value = "Territoriality (females)"
$ws.Range("F1").Value = "Group Size During Reproduction"
$ws.Range("G1").Value = "Group Size Outside of Reproduction"
$ws.Range("H1").Value = "Group Property"
$ws.Range("I1").Value = "Age at Maturity"
$ws.Range("J1").Value = "Avg. life expectancy"
$ws.Range("K1").Value = "# offspring/reproductive bout"
$ws.Range("L1").Value = "# reproductive bouts/year"
$ws.Range("M1").Value = "Migratory behavior"
$ws.Range("N1").Value = "Activity Pattern"
$ws.Range("O1").Value = "Habitat Complexity"

# --- Set the species column (A2:A24) values ---
$ws.Range("A2").Value = "Microtus ochrogaster"
$ws.Range("A3").Value = "Microtus pennsylvanicus"
$ws.Range("A4").Value = "Mus musculus"
$ws.Range("A5").Value = "Rattus norvegicus"
$ws.Range("A6").Value = "Cavia porcellus"
$ws.Range("A7").Value = "Fukomys damarensis"
$ws.Range("A8").Value = "Pan troglodytes"
$ws.Range("A9").Value = "Homo sapiens"
$ws.Range("A10").Value = "Gorilla gorilla"
$ws.Range("A11").Value = "Hylobates lar"
$ws.Range("A12").Value = "Macaca mulatta"
$ws.Range("A13").Value = "Macaca fascicularis"
$ws.Range("A14").Value = "Passer domesticus"
$ws.Range("A15").Value = "Zonotrichia albicollis"
$ws.Range("A16").Value = "Taeniopygia guttata"
$ws.Range("A17").Value = "Pipra filicauda"
$ws.Range("A18").Value = "Columba livia"
$ws.Range("A19").Value = "Coturnix japonica"
$ws.Range("A20").Value = "Anolis carolinensis"
$ws.Range("A21").Value = "Engystomops putulosus"
$ws.Range("A22").Value = "Dendrobate auratus"
$ws.Range("A23").Value = "Astatotilap burtoni"
$ws.Range("A24").Value = "Gasterosteus aculeatus"
